$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Ref, $NewValue)
    # Force the cell to be stored as literal text even when the
    # string looks numeric (e.g. "7.05"), matching the source file
    # which keeps these Price values as inline strings, not numbers.
    $Sheet.Range($Ref).NumberFormat = "@"
    $Sheet.Range($Ref).Value = $NewValue
    $Sheet.Range($Ref).Style = "Normal"
}

Set-TextValue $ws 'D2' '67.840.68'
$ws.Range('E2').Value = '  +1.10%  '
Set-TextValue $ws 'D3' '3.246.74'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws 'D5' '582.48'
$ws.Range('E5').Value = '  +0.85%  '
Set-TextValue $ws 'D6' '182.91'
$ws.Range('E6').Value = '  +3.80%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('E9').Value = '  +4.21%  '
Set-TextValue $ws 'D10' '6.68'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +1.92%  '
Set-TextValue $ws 'D12' '3.808.69'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +0.45%  '
Set-TextValue $ws 'D14' '28.71'
$ws.Range('E14').Value = '  +3.24%  '
Set-TextValue $ws 'D15' '67.845.70'
$ws.Range('E15').Value = '  +1.20%  '
Set-TextValue $ws 'D16' '0.0000171'
$ws.Range('E16').Value = '  +2.26%  '
Set-TextValue $ws 'D17' '3.251.26'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('E19').Value = '  +2.04%  '
Set-TextValue $ws 'D20' '379.49'
$ws.Range('E20').Value = '  +3.27%  '
$ws.Range('E21').Value = '  +2.50%  '
Set-TextValue $ws 'D23' '71.35'
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('E25').Value = '  +0.65%  '
Set-TextValue $ws 'D26' '9.91'
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('E27').Value = '  +2.23%  '
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('E29').Value = '  +0.44%  '
Set-TextValue $ws 'D30' '5.66'
$ws.Range('E30').Value = '  +0.80%  '
Set-TextValue $ws 'D31' '22.85'
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D32' '7.05'
$ws.Range('E32').Value = '  +4.34%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws 'D33' '0.998'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  +2.57%  '
$ws.Range('E35').Value = '  +4.24%  '
Set-TextValue $ws 'D36' '162.24'
$ws.Range('E36').Value = '  -5.74%  '
Set-TextValue $ws 'D37' '0.839'
$ws.Range('E37').Value = '  -1.19%  '
Set-TextValue $ws 'D38' '1.85'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D39' '26.45'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D40' '6.70'
$ws.Range('E40').Value = '  +5.28%  '
$ws.Range('E41').Value = '  +6.87%  '
$ws.Range('E42').Value = '  +0.25%  '
Set-TextValue $ws 'D43' '25.48'
$ws.Range('E43').Value = '  +3.85%  '
Set-TextValue $ws 'D44' '41.17'
$ws.Range('E44').Value = '  +1.98%  '
Set-TextValue $ws 'D45' '345.80'
$ws.Range('E45').Value = '  +3.88%  '
Set-TextValue $ws 'D46' '0.0684'
$ws.Range('E46').Value = '  +2.02%  '
Set-TextValue $ws 'D47' '2.626.68'
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('E48').Value = '  +1.42%  '
$ws.Range('E49').Value = '  -0.98%  '
Set-TextValue $ws 'D50' '0.993'
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('E51').Value = '  +2.52%  '
